$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate data among rows 2, 3, 5 (D, J, K, M, P columns) per commit diff:
# new Row2 = old Row3, new Row3 = old Row5, new Row5 = old Row2

$row2 = @{ D = $ws.Range("D2").Value2; J = $ws.Range("J2").Value2; K = $ws.Range("K2").Value2; M = $ws.Range("M2").Value2; P = $ws.Range("P2").Value2 }
$row3 = @{ D = $ws.Range("D3").Value2; J = $ws.Range("J3").Value2; K = $ws.Range("K3").Value2; M = $ws.Range("M3").Value2; P = $ws.Range("P3").Value2 }
$row5 = @{ D = $ws.Range("D5").Value2; J = $ws.Range("J5").Value2; K = $ws.Range("K5").Value2; M = $ws.Range("M5").Value2; P = $ws.Range("P5").Value2 }

$ws.Range("D2").Value2 = $row3.D
$ws.Range("J2").Value2 = $row3.J
$ws.Range("K2").Value2 = $row3.K
$ws.Range("M2").Value2 = $row3.M
$ws.Range("P2").Value2 = $row3.P

$ws.Range("D3").Value2 = $row5.D
$ws.Range("J3").Value2 = $row5.J
$ws.Range("K3").Value2 = $row5.K
$ws.Range("M3").Value2 = $row5.M
$ws.Range("P3").Value2 = $row5.P

$ws.Range("D5").Value2 = $row2.D
$ws.Range("J5").Value2 = $row2.J
$ws.Range("K5").Value2 = $row2.K
$ws.Range("M5").Value2 = $row2.M
$ws.Range("P5").Value2 = $row2.P
